# "Add files via upload" - append the next batch of tracked properties to
# the bottom of the Sheet1 tracking list (rows 145-156), then clear the
# now-superfluous fill style that used to be applied to every B-column
# cell, and leave the selection on the last-entered cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New rows (ID in column A, 단지명/name in column B) -------------------
$newRows = @(
    @(110209, "e편한세상신촌"),
    @(110092, "신촌푸르지오"),
    @(849,    "두산"),
    @(26046,  "돈의문센트레빌"),
    @(11567,  "천연뜨란채"),
    @(941,    "독립문극동"),
    @(108064, "DMC파크뷰자이"),
    @(124802, "홍제역해링턴플레이스"),
    @(945,    "인왕산현대"),
    @(128515, "e편한세상서대문"),
    @(114768, "연희파크푸르지오"),
    @(128027, "힐스테이트홍은포레스트")
)

$startRow = 145
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $id = $newRows[$i][0]
    $name = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $name
}

# --- Drop the stray applyFill style that every old B-column cell carried --
$oldRange = $ws.Range("B2:B144")
$oldRange.ClearFormats()

# --- Leave the selection/scroll where data entry finished ------------------
$ws.Range("B153").Select()

Write-Host "done"
